$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Resolve item: "OX" -> "oxidized residue" for the modification type in B2
$ws.Range("B2").Value = "oxidized residue"

# Reflect the new active selection recorded in the saved view state
[void]$ws.Range("B5").Select()
